# Fixed naive component forecaster bug - Presentation state 11.02.
# Updates the ME/MAE/MSE/RMSE/SE/N values for rows 2-11 (Q0..Q9) on sheet T45.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.1352154364003263, 0.23198446106715, 0.2545140212335539, 0.504493826754653, 0.4969593252438608, 23)
    3  = @(0.410665070665073, 2.596388298586061, 26.31885526302818, 5.130190567905658, 5.234067016424154, 22)
    4  = @(-0.9427454607211139, 2.394021067790873, 19.00138225390311, 4.359057496053833, 4.360991161731755, 21)
    5  = @(-0.1609758476115556, 1.331657534217984, 6.400840711354078, 2.529988282849168, 2.590453625757429, 20)
    6  = @(-0.07574193521395585, 1.494237072137884, 9.671523505559456, 3.109907314625222, 3.194178268796918, 19)
    7  = @(-0.3467465960030111, 1.917003265407361, 13.48431232460344, 3.672099171400936, 3.761675348046774, 18)
    8  = @(-0.3096017392693358, 1.877145454039581, 10.57116444722301, 3.251332718628318, 3.336168185344669, 17)
    9  = @(-0.0424536219110716, 1.743620181688103, 11.37086449867186, 3.372071247567562, 3.482384193035339, 16)
    10 = @(-0.1474229548386492, 1.844192231578261, 11.76371683194192, 3.429827522185616, 3.546927740664431, 15)
    11 = @(-0.1432630465947277, 2.238041655440814, 14.73581901823137, 3.838726223401633, 3.992689479225886, 13)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
